$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet (sheet 1) ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.1.0 -> 2.2.0-ballot
$meta.Range("B3").Value = "2.2.0-ballot"

# Date: 2025-12-18T17:25:31+00:00 -> 2025-12-19T08:32:44+00:00
$meta.Range("B8").Value = "2025-12-19T08:32:44+00:00"

# Base Definition: append FHIR version pin
$meta.Range("B18").Value = "http://hl7.org/fhir/StructureDefinition/Extension|4.0.1"

# --- "Elements" sheet (sheet 2) ---
$elements = $wb.Worksheets.Item("Elements")

# Binding Value Set (Z6): append value set version pin
$elements.Range("Z6").Value = "https://smt.esante.gouv.fr/fhir/ValueSet/jdv-motif-non-realisation-evenement-cisis|20250624152100"

# Column Z widened (now bestFit's to a longer value) from 65.457.. to 79.988..
$elements.Columns.Item(26).ColumnWidth = 79.2
